$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Monthly update (MV -datos-): correct the last two existing rows and
# --- append the new row for 01-07-2021 ---

# Row 174 -> Serie 01-05-2021 : revised figures
$ws.Range("B174").Value = 63044
$ws.Range("D174").Value = 53995
$ws.Range("E174").Value = 14550
$ws.Range("G174").Value = 12665

# Row 175 -> Serie 01-06-2021 : revised figures
$ws.Range("B175").Value = 59033
$ws.Range("D175").Value = 49938

# Row 176 -> new Serie 01-07-2021
# Force the date-like label to be stored as text (matches the existing
# "Serie" column, which holds shared strings, not real dates), then put the
# cell style back the way it was on the previous row so no visible
# formatting change is introduced.
$ws.Range("A176").NumberFormat = "@"
$ws.Range("A176").Value = "01-07-2021"
$ws.Range("A176").Style = $ws.Range("A175").Style

$ws.Range("B176").Value = 56021
$ws.Range("C176").Value = 8875
$ws.Range("D176").Value = 47146
$ws.Range("E176").Value = 10739
$ws.Range("F176").Value = 2112
$ws.Range("G176").Value = 8627
